$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Player")

# Update Attack (column C) values per new stat balance
$ws.Range("C2").Value = 40
$ws.Range("C3").Value = 80
$ws.Range("C4").Value = 120
$ws.Range("C6").Value = 200
$ws.Range("C7").Value = 240
$ws.Range("C8").Value = 280
$ws.Range("C9").Value = 320
$ws.Range("C10").Value = 360
$ws.Range("C11").Value = 400

# Match the saved selection state (C2:C11 selected, active cell C2)
$ws.Range("C2:C11").Select()
